$wb = $excel.ActiveWorkbook

# --- Sheet: Change Management Overview ---
$ws1 = $wb.Worksheets.Item("Change Management Overview")

# Project name revert
$ws1.Range("B6").Value = "Enterprise Cloud Infrastructure Migration"

# Restore blank separator rows (no attributes, no shift of following rows)
$ws1.Rows(13).OutlineLevel = 0
$ws1.Rows(21).OutlineLevel = 0

# Objective text reverts (AI/ML -> IT)
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new IT systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in IT technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for IT transformation"

# --- Sheet: Change Impact Assessment ---
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

$ws2.Rows(2).OutlineLevel = 0

$ws2.Range("A4").Value = "IT Managers"
$ws2.Range("G4").Value = "IT automation"
$ws2.Range("A5").Value = "System Administrators"

# --- Sheet: Change Activities ---
$ws3 = $wb.Worksheets.Item("Change Activities")

$ws3.Rows(2).OutlineLevel = 0
